$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rallies")

# Update existing row 72: rally_no (column D) goes from 9 to 10
$ws.Cells.Item(72, 4).Value = 10

# Append a new row 73, duplicating row 72's content but advancing the counters
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = 1
$ws.Cells.Item(73, 3).Value = 3
$ws.Cells.Item(73, 4).Value = 10
$ws.Cells.Item(73, 5).Value = "NOS"
$ws.Cells.Item(73, 6).Value = ""
$ws.Cells.Item(73, 7).Value = 4
$ws.Cells.Item(73, 8).Value = "MEIO"
$ws.Cells.Item(73, 9).Value = "PONTO"
$ws.Cells.Item(73, 10).Value = "NOS"
$ws.Cells.Item(73, 11).Value = 10
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = "1 4 m"
$ws.Cells.Item(73, 14).Value = "FRENTE"
$ws.Cells.Item(73, 15).Value = "FRENTE"
$ws.Cells.Item(73, 16).Value = "FRENTE"
